$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:G190")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("F2:F190"), 0, 2) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("G2:G190"), 0, 2) | Out-Null

$ws.Sort.SetRange($rng)
$ws.Sort.Header = -4142  # xlNo
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1  # xlTopToBottom
$ws.Sort.Apply()
